$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet rename: "PCBA BOM Template" -> "POWER PCBA" -----------------
# Renaming the sheet updates <sheet name="..."> and the _FilterDatabase
# defined name automatically. The Print_Area defined name is re-pointed
# explicitly afterwards (re-asserting PageSetup.PrintArea also refreshes
# the sheet-qualifier used by the _xlnm.Print_Area defined name).
$ws.Name = "POWER PCBA"
$ws.PageSetup.PrintArea = "A1:N15"

# --- Selection moves from A1:N11 (active cell N11) to D28 --------------
$ws.Range("D28").Select()

# --- Shared formula N7 (si=0) no longer spans down to row 13 -----------
# Originally the shared formula group covering N7:N13 only had live
# formulas in N7:N9 (N10:N13 were already blank); clear any stray
# content in that tail so just N7:N9 remain populated.
$ws.Range("N10:N13").ClearContents()
